# Update cryptocurrency price/volume values (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.272.37'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '1.651.67'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  -0.49%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.515'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.10%  '
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.257'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.02%  '
$ws.Range('E9').Value = '  +1.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0849'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.22%  '
$ws.Range('D12').Value = '1.881.80'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').Value = '1.659.09'
$ws.Range('E13').Value = '  +0.47%  '
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('E15').Value = '  +2.63%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.75'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').Value = '27.252.99'
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '220.38'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.86'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.55'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.14%  '
$ws.Range('E23').Value = '  +0.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.23'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.98'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.55'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.23%  '
$ws.Range('E27').Value = '  -0.38%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('E29').Value = '  -0.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0510'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.37%  '
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('E33').Value = '  +1.70%  '
$ws.Range('E34').Value = '  +1.61%  '
$ws.Range('D35').Value = '1.268.90'
$ws.Range('E35').Value = '  +1.02%  '
$ws.Range('E37').Value = '  +0.95%  '
$ws.Range('E38').Value = '  +2.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.846'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.44%  '
$ws.Range('E40').Value = '  -0.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.811'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.48%  '
$ws.Range('E42').Value = '  +1.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.21'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.14%  '
$ws.Range('D44').Value = '1.791.96'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.18'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.89'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.61'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('D48').Value = '0.0₆0106'
$ws.Range('E48').Value = '  +7.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0514'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.69'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0972'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.44%  '
